$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header fields
$ws.Range("B2").Value = "Repair HVAC Interior B10210"

# F2 holds a literal text date string ("08/22/2022") even though the cell is
# number-formatted as a date. A plain .Value assignment would get silently
# parsed into a real date serial (wrong type), so stage the text in a
# scratch cell forced to Text format, then copy only the value across -
# this preserves F2's original style/format and keeps it a text cell.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "08/22/2022"
$ws.Range("Z1").Copy()
$ws.Range("F2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
# PasteSpecial into a merged cell drops the merge - restore it.
$ws.Range("F2:G2").Merge()

$ws.Range("B3").Value = "RKMF 17-0133"
$ws.Range("F3").Value = 5733
$ws.Range("B4").Value = "FA486118FA140"

# Manpower row 7
$ws.Range("A7").Value = "548 Group Inc."
$ws.Range("B7").Value = "Bldg 400"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "  Abatement"

# Clear "No on-site work" note in E15
$ws.Range("E15").Value = ""

# Inspections row 18
$ws.Range("A18").Value = "Inpsector A"
$ws.Range("B18").Value = "LAISD"
$ws.Range("C18").Value = "Bldg 400"
$ws.Range("E18").Value = "Rough Electrical "
$ws.Range("G18").Value = "Pass"

# Note
$ws.Range("A22").Value = "Test"
